# Added filtering options for the Component Analysis
#
# A new evaluation row (Q0) is inserted at the top of the data block (row 2),
# which shifts all existing Q0..Q8 rows down by one row (now Q1..Q9), and the
# previously last row (old Q9, row 11) is dropped from the table. The Q-labels
# in column A stay fixed in place; only the statistic columns B:G move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 10 -> 11, 9 -> 10, ... 2 -> 3 (process bottom-up so
# we never overwrite a row before it has been copied down).
for ($r = 10; $r -ge 2; $r--) {
    $src = $ws.Range("B$r`:G$r")
    $dst = $ws.Range("B$($r+1)`:G$($r+1)")
    $dst.Value2 = $src.Value2
}

# Fill in the new first row (Q0) with the newly computed statistics.
$ws.Range("B2").Value = 0.03483647684766893
$ws.Range("C2").Value = 0.9165303275553447
$ws.Range("D2").Value = 4.351217448857517
$ws.Range("E2").Value = 2.085957202067558
$ws.Range("F2").Value = 2.108713024898596
$ws.Range("G2").Value = 46
